$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2:N5").Value = 53.71147335634279
